$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.01
$ws.Range("F2").Value = 0.0011
$ws.Range("I2").Value = 12000

$ws.Range("I3").Value = 100

$ws.Range("I4").Value = 30100

$ws.Range("G5").Value = -0.64
$ws.Range("I5").Value = 2000

$ws.Range("E6").Value = 0.02
$ws.Range("F6").Value = 0.0019
$ws.Range("I6").Value = 200

$ws.Range("I8").Value = 3100

$ws.Range("I9").Value = 26700

$ws.Range("G10").Value = -0.55
$ws.Range("I10").Value = 20100

$ws.Range("I11").Value = 280100

$ws.Range("E12").Value = 0.03
$ws.Range("F12").Value = 0.0022
$ws.Range("I12").Value = 71700

$ws.Range("I13").Value = 7500

$ws.Range("E14").Value = 0.03
$ws.Range("F14").Value = 0.0021
$ws.Range("I14").Value = 35500

$ws.Range("G15").Value = -0.46
$ws.Range("I15").Value = 7800

$ws.Range("G16").Value = -0.45

$ws.Range("G17").Value = -0.44
$ws.Range("H17").Value = 0.01

$ws.Range("E18").Value = 0.04
$ws.Range("F18").Value = 0.0025
$ws.Range("G18").Value = -0.43
$ws.Range("I18").Value = 1800

$ws.Range("E19").Value = 0.04
$ws.Range("F19").Value = 0.0025
$ws.Range("I19").Value = 7900

$ws.Range("E20").Value = 0.06
$ws.Range("F20").Value = 0.0037
$ws.Range("H20").Value = 0.04
$ws.Range("I20").Value = 10000

$ws.Range("E21").Value = 0.04
$ws.Range("F21").Value = 0.0024
$ws.Range("H21").Value = 0.06
$ws.Range("I21").Value = 400

$ws.Range("E22").Value = 0.06
$ws.Range("F22").Value = 0.0036
$ws.Range("H22").Value = 0.08
$ws.Range("I22").Value = 500

$ws.Range("E23").Value = 0.07000000000000001
$ws.Range("F23").Value = 0.0041
$ws.Range("H23").Value = 0.11
$ws.Range("I23").Value = 2300

$ws.Range("H24").Value = 0.16

$ws.Range("E25").Value = 0.06
$ws.Range("F25").Value = 0.0034
$ws.Range("H25").Value = 0.21
$ws.Range("I25").Value = 122100

$ws.Range("G26").Value = -0.35
$ws.Range("H26").Value = 0.37
$ws.Range("I26").Value = 48300

$ws.Range("E27").Value = 0.08
$ws.Range("F27").Value = 0.0044
$ws.Range("G27").Value = -0.34
$ws.Range("H27").Value = 0.48
$ws.Range("I27").Value = 84000

$ws.Range("E28").Value = 0.1
$ws.Range("F28").Value = 0.0054
$ws.Range("G28").Value = -0.33
$ws.Range("H28").Value = 0.61
$ws.Range("I28").Value = 1600

$ws.Range("H29").Value = 0.78

$ws.Range("E30").Value = 0.11
$ws.Range("F30").Value = 0.0058
$ws.Range("H30").Value = 0.98
$ws.Range("I30").Value = 54600

$ws.Range("H31").Value = 1.23
$ws.Range("I31").Value = 103600

$ws.Range("H32").Value = 1.52
$ws.Range("I32").Value = 328100

$ws.Range("E33").Value = 0.13
$ws.Range("F33").Value = 0.0066
$ws.Range("H33").Value = 1.87
$ws.Range("I33").Value = 69700

$ws.Range("E34").Value = 0.13
$ws.Range("F34").Value = 0.0065
$ws.Range("H34").Value = 2.27
$ws.Range("I34").Value = 249300

$ws.Range("H35").Value = 2.74
$ws.Range("I35").Value = 163500

$ws.Range("E36").Value = 0.15
$ws.Range("F36").Value = 0.0073
$ws.Range("G36").Value = -0.26
$ws.Range("H36").Value = 3.28
$ws.Range("I36").Value = 60900

$ws.Range("G37").Value = -0.25
$ws.Range("H37").Value = 3.89
$ws.Range("I37").Value = 22200

$ws.Range("E38").Value = 0.22
$ws.Range("F38").Value = 0.0105
$ws.Range("G38").Value = -0.24
$ws.Range("H38").Value = 4.59
$ws.Range("I38").Value = 182300

$ws.Range("G39").Value = -0.23
$ws.Range("H39").Value = 5.37
$ws.Range("I39").Value = 12900

$ws.Range("G40").Value = -0.22
$ws.Range("H40").Value = 6.25
$ws.Range("I40").Value = 96700

$ws.Range("E41").Value = 0.22
$ws.Range("F41").Value = 0.0101
$ws.Range("H41").Value = 7.22
$ws.Range("I41").Value = 59600

$ws.Range("E42").Value = 0.24
$ws.Range("F42").Value = 0.0109
$ws.Range("H42").Value = 8.289999999999999
$ws.Range("I42").Value = 360100

$ws.Range("H43").Value = 9.460000000000001
$ws.Range("I43").Value = 625500

$ws.Range("E44").Value = 0.28
$ws.Range("F44").Value = 0.0124
$ws.Range("H44").Value = 10.73
$ws.Range("I44").Value = 286300

$ws.Range("E45").Value = 0.3
$ws.Range("F45").Value = 0.0132
$ws.Range("H45").Value = 12.11
$ws.Range("I45").Value = 217900

$ws.Range("E46").Value = 0.38
$ws.Range("F46").Value = 0.0165
$ws.Range("H46").Value = 13.59
$ws.Range("I46").Value = 361500

$ws.Range("E47").Value = 0.36
$ws.Range("F47").Value = 0.0155
$ws.Range("G47").Value = -0.16
$ws.Range("H47").Value = 15.17
$ws.Range("I47").Value = 269100

$ws.Range("E48").Value = 0.41
$ws.Range("F48").Value = 0.0174
$ws.Range("G48").Value = -0.15
$ws.Range("H48").Value = 16.85
$ws.Range("I48").Value = 120400

$ws.Range("E49").Value = 0.42
$ws.Range("F49").Value = 0.0177
$ws.Range("G49").Value = -0.14
$ws.Range("H49").Value = 18.62
$ws.Range("I49").Value = 203300

$ws.Range("D50").Value = "Montar"
$ws.Range("E50").Value = 0.5
$ws.Range("F50").Value = 0.0208
$ws.Range("G50").Value = -0.13
$ws.Range("H50").Value = 20.49
$ws.Range("I50").Value = 351100

$ws.Range("E51").Value = 0.55
$ws.Range("F51").Value = 0.0227
$ws.Range("G51").Value = -0.12
$ws.Range("H51").Value = 22.44
$ws.Range("I51").Value = 101100

$ws.Range("D52").Value = "Montar"
$ws.Range("E52").Value = 0.66
$ws.Range("F52").Value = 0.0267
$ws.Range("H52").Value = 26.58
$ws.Range("I52").Value = 132300

$ws.Range("E53").Value = 0.73
$ws.Range("F53").Value = 0.0292
$ws.Range("H53").Value = 28.75
$ws.Range("I53").Value = 460400

$ws.Range("E54").Value = 0.77
$ws.Range("F54").Value = 0.0305
$ws.Range("H54").Value = 30.98
$ws.Range("I54").Value = 127600

$ws.Range("H55").Value = 33.25
$ws.Range("I55").Value = 246600

$ws.Range("H56").Value = 33.25

$ws.Range("E57").Value = 0.9
$ws.Range("F57").Value = 0.035
$ws.Range("H57").Value = 35.57
$ws.Range("I57").Value = 360100

$ws.Range("E58").Value = 1.03
$ws.Range("F58").Value = 0.0396
$ws.Range("G58").Value = -0.06
$ws.Range("H58").Value = 37.92
$ws.Range("I58").Value = 1527600

$ws.Range("E59").Value = 1.1
$ws.Range("F59").Value = 0.0419
$ws.Range("G59").Value = -0.05
$ws.Range("H59").Value = 40.29
$ws.Range("I59").Value = 216600

$ws.Range("G60").Value = -0.04
$ws.Range("H60").Value = 42.68
$ws.Range("I60").Value = 232000

$ws.Range("E61").Value = 1.33
$ws.Range("F61").Value = 0.0497
$ws.Range("G61").Value = -0.03
$ws.Range("H61").Value = 45.07
$ws.Range("I61").Value = 277200

$ws.Range("E62").Value = 1.35
$ws.Range("F62").Value = 0.05
$ws.Range("G62").Value = -0.02
$ws.Range("H62").Value = 47.45
$ws.Range("I62").Value = 730400

$ws.Range("E63").Value = 1.67
$ws.Range("F63").Value = 0.0607
$ws.Range("G63").Value = -0.0
$ws.Range("H63").Value = 52.17
$ws.Range("I63").Value = 809800

$ws.Range("G64").Value = -0.0
$ws.Range("H64").Value = 52.17

$ws.Range("E65").Value = 1.8
$ws.Range("F65").Value = 0.0649
$ws.Range("H65").Value = 54.5
$ws.Range("I65").Value = 230100

$ws.Range("E66").Value = 1.93
$ws.Range("F66").Value = 0.0689
$ws.Range("H66").Value = 56.78
$ws.Range("I66").Value = 1150500

$ws.Range("E67").Value = 2.23
$ws.Range("F67").Value = 0.07820000000000001
$ws.Range("H67").Value = 61.22
$ws.Range("I67").Value = 148200

$ws.Range("H68").Value = 61.22

$ws.Range("E69").Value = 2.34
$ws.Range("F69").Value = 0.0814
$ws.Range("G69").Value = 0.04
$ws.Range("H69").Value = 63.36
$ws.Range("I69").Value = 104700

$ws.Range("E70").Value = 2.58
$ws.Range("F70").Value = 0.089
$ws.Range("G70").Value = 0.05
$ws.Range("H70").Value = 65.44
$ws.Range("I70").Value = 1241600

$ws.Range("E71").Value = 2.8
$ws.Range("F71").Value = 0.0949
$ws.Range("G71").Value = 0.07000000000000001
$ws.Range("H71").Value = 69.42
$ws.Range("I71").Value = 67200

$ws.Range("E72").Value = 3.24
$ws.Range("F72").Value = 0.108
$ws.Range("G72").Value = 0.09
$ws.Range("H72").Value = 73.11
$ws.Range("I72").Value = 95000

$ws.Range("E73").Value = 3.62
$ws.Range("F73").Value = 0.1187
$ws.Range("H73").Value = 76.51000000000001
$ws.Range("I73").Value = 300500

$ws.Range("H74").Value = 76.51000000000001

$ws.Range("E75").Value = 3.43
$ws.Range("F75").Value = 0.1115
$ws.Range("H75").Value = 78.09999999999999
$ws.Range("I75").Value = 2000

$ws.Range("E76").Value = 3.91
$ws.Range("F76").Value = 0.1261
$ws.Range("H76").Value = 79.61
$ws.Range("I76").Value = 3800

$ws.Range("G77").Value = 0.14
$ws.Range("H77").Value = 82.41

$ws.Range("E78").Value = 4.61
$ws.Range("F78").Value = 0.1452
$ws.Range("G78").Value = 0.15
$ws.Range("H78").Value = 83.7
$ws.Range("I78").Value = 146800

$ws.Range("G79").Value = 0.16
$ws.Range("H79").Value = 84.91

$ws.Range("G80").Value = 0.18
$ws.Range("H80").Value = 87.14

$ws.Range("E81").Value = 5.78
$ws.Range("F81").Value = 0.1752
$ws.Range("G81").Value = 0.2
$ws.Range("H81").Value = 89.09
$ws.Range("I81").Value = 71600
